$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking strings (e.g. "0.656")
# round-trip exactly instead of being coerced to floating point numbers.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '36.130.24'
$ws.Range('E2').Value = '  -1.60%  '
$ws.Range('D3').Value = '2.022.92'
$ws.Range('E3').Value = '  -2.74%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').Value = '243.38'
$ws.Range('E5').Value = '  -0.61%  '
$ws.Range('D6').Value = '0.656'
$ws.Range('E6').Value = '  +0.79%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').Value = '54.32'
$ws.Range('E8').Value = '  +0.47%  '
$ws.Range('D9').Value = '58.77'
$ws.Range('E9').Value = '  +0.18%  '
$ws.Range('D10').Value = '0.361'
$ws.Range('E10').Value = '  -1.17%  '
$ws.Range('D11').Value = '0.0734'
$ws.Range('E11').Value = '  -3.45%  '
$ws.Range('E12').Value = '  -4.19%  '
$ws.Range('D13').Value = '0.892'
$ws.Range('E13').Value = '  +0.11%  '
$ws.Range('D14').Value = '14.04'
$ws.Range('E14').Value = '  -6.18%  '
$ws.Range('D15').Value = '2.324.02'
$ws.Range('E15').Value = '  -2.59%  '
$ws.Range('D16').Value = '5.29'
$ws.Range('E16').Value = '  -3.47%  '
$ws.Range('D17').Value = '2.043.70'
$ws.Range('E17').Value = '  -1.27%  '
$ws.Range('D18').Value = '17.30'
$ws.Range('E18').Value = '  +0.57%  '
$ws.Range('D19').Value = '35.998.33'
$ws.Range('E19').Value = '  -1.88%  '
$ws.Range('D20').Value = '71.04'
$ws.Range('E20').Value = '  -1.97%  '
$ws.Range('D21').Value = '0.0₃0846'
$ws.Range('E21').Value = '  -3.55%  '
$ws.Range('D22').Value = '235.53'
$ws.Range('E22').Value = '  -1.46%  '
$ws.Range('D23').Value = '5.13'
$ws.Range('E23').Value = '  -5.42%  '
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('E25').Value = '  -1.93%  '
$ws.Range('D26').Value = '2.28'
$ws.Range('E26').Value = '  +5.72%  '
$ws.Range('E27').Value = '  -6.47%  '
$ws.Range('D28').Value = '162.99'
$ws.Range('E28').Value = '  -2.77%  '
$ws.Range('D29').Value = '19.74'
$ws.Range('E29').Value = '  -4.10%  '
$ws.Range('D30').Value = '0.120'
$ws.Range('E30').Value = '  -2.25%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = '4.90'
$ws.Range('E31').Value = '  -7.58%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').Value = '1.15'
$ws.Range('E32').Value = '  -2.56%  '
$ws.Range('D33').Value = '0.0593'
$ws.Range('E33').Value = '  -1.99%  '
$ws.Range('D34').Value = '4.31'
$ws.Range('E34').Value = '  -7.47%  '
$ws.Range('D35').Value = '0.0886'
$ws.Range('E35').Value = '  +7.60%  '
$ws.Range('E36').Value = '  +0.12%  '
$ws.Range('E37').Value = '  -1.36%  '
$ws.Range('D38').Value = '2.17'
$ws.Range('E38').Value = '  -8.65%  '
$ws.Range('D39').Value = '5.00'
$ws.Range('E39').Value = '  +3.85%  '
$ws.Range('D40').Value = '1.19'
$ws.Range('E40').Value = '  -5.09%  '
$ws.Range('D41').Value = '2.89'
$ws.Range('E41').Value = '  +1.48%  '
$ws.Range('D42').Value = '0.0213'
$ws.Range('E42').Value = '  -3.08%  '
$ws.Range('E43').Value = '  -5.37%  '
$ws.Range('D44').Value = '0.0895'
$ws.Range('E44').Value = '  -5.86%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = '91.47'
$ws.Range('E45').Value = '  -4.56%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = '1.384.13'
$ws.Range('E46').Value = '  +0.86%  '
$ws.Range('D47').Value = '7.38'
$ws.Range('E47').Value = '  +0.95%  '
$ws.Range('D48').Value = '15.35'
$ws.Range('E48').Value = '  -3.85%  '
$ws.Range('D49').Value = '2.93'
$ws.Range('E49').Value = '  +0.84%  '
$ws.Range('D50').Value = '2.23'
$ws.Range('E50').Value = '  -8.02%  '
$ws.Range('D51').Value = '45.32'
$ws.Range('E51').Value = '  +0.67%  '

# Restore the default (Normal) style so cell formatting matches the original
# workbook - only the text content changed, not the number format.
$ws.Range('D2:D51').Style = 'Normal'
